$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "70.316.88"
$ws.Range("E2").Value = "  +0.42%  "

# Row 3
$ws.Range("D3").Value = "3.609.40"
$ws.Range("E3").Value = "  +1.60%  "

# Row 4
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.07%  "

# Row 5
$ws.Range("D5").Value = "'605.22"
$ws.Range("E5").Value = "  +0.25%  "

# Row 6
$ws.Range("D6").Value = "'196.24"
$ws.Range("E6").Value = "  -0.26%  "

# Row 7
$ws.Range("D7").Value = "'0.627"
$ws.Range("E7").Value = "  -0.02%  "

# Row 8
$ws.Range("D8").Value = "'1.00"

# Row 9
$ws.Range("D9").Value = "'0.207"
$ws.Range("E9").Value = "  -1.22%  "

# Row 10
$ws.Range("D10").Value = "'0.646"
$ws.Range("E10").Value = "  -1.38%  "

# Row 11
$ws.Range("D11").Value = "'53.75"
$ws.Range("E11").Value = "  -0.75%  "

# Row 12
$ws.Range("D12").Value = "'0.0000304"
$ws.Range("E12").Value = "  +0.16%  "

# Row 13
$ws.Range("D13").Value = "'9.56"
$ws.Range("E13").Value = "  -0.05%  "

# Row 14
$ws.Range("D14").Value = "4.178.93"
$ws.Range("E14").Value = "  +1.47%  "

# Row 15
$ws.Range("D15").Value = "'13.04"
$ws.Range("E15").Value = "  +2.68%  "

# Row 16
$ws.Range("D16").Value = "'595.93"
$ws.Range("E16").Value = "  -1.03%  "

# Row 17
$ws.Range("D17").Value = "70.491.73"
$ws.Range("E17").Value = "  +0.41%  "

# Row 18
$ws.Range("D18").Value = "3.605.33"
$ws.Range("E18").Value = "  +1.49%  "

# Row 19
$ws.Range("D19").Value = "'19.03"
$ws.Range("E19").Value = "  -0.67%  "

# Row 20
$ws.Range("D20").Value = "'0.123"
$ws.Range("E20").Value = "  +1.47%  "

# Row 21
$ws.Range("D21").Value = "'0.998"
$ws.Range("E21").Value = "  -0.03%  "

# Row 22
$ws.Range("D22").Value = "'17.84"
$ws.Range("E22").Value = "  -1.08%  "

# Row 23
$ws.Range("E23").Value = "  -1.75%  "

# Row 24
$ws.Range("D24").Value = "'102.24"
$ws.Range("E24").Value = "  -1.01%  "

# Row 25
$ws.Range("D25").Value = "'4.63"
$ws.Range("E25").Value = "  -0.18%  "

# Row 26
$ws.Range("E26").Value = "  -3.09%  "

# Row 27
$ws.Range("D27").Value = "'10.78"
$ws.Range("E27").Value = "  -1.63%  "

# Row 28
$ws.Range("D28").Value = "'9.64"
$ws.Range("E28").Value = "  -0.50%  "

# Row 29
$ws.Range("D29").Value = "'33.89"
$ws.Range("E29").Value = "  +0.30%  "

# Row 30
$ws.Range("D30").Value = "'4.74"
$ws.Range("E30").Value = "  +4.32%  "

# Row 31
$ws.Range("D31").Value = "'7.22"
$ws.Range("E31").Value = "  +1.26%  "

# Row 32
$ws.Range("D32").Value = "'12.30"
$ws.Range("E32").Value = "  -3.39%  "

# Row 33
$ws.Range("D33").Value = "'0.118"
$ws.Range("E33").Value = "  +1.54%  "

# Row 34
$ws.Range("B34").Value = "PEPE"
$ws.Range("C34").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D34").Value = "0.0₃0903"
$ws.Range("E34").Value = "  +4.06%  "

# Row 35
$ws.Range("B35").Value = "OKB"
$ws.Range("C35").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D35").Value = "'63.19"
$ws.Range("E35").Value = "  -0.42%  "

# Row 36
$ws.Range("D36").Value = "3.891.66"
$ws.Range("E36").Value = "  +3.66%  "

# Row 37
$ws.Range("D37").Value = "'538.19"
$ws.Range("E37").Value = "  +9.12%  "

# Row 38
$ws.Range("D38").Value = "'3.12"
$ws.Range("E38").Value = "  +0.95%  "

# Row 39
$ws.Range("E39").Value = "  -0.13%  "

# Row 40
$ws.Range("D40").Value = "'37.02"
$ws.Range("E40").Value = "  +0.29%  "

# Row 41
$ws.Range("D41").Value = "'0.392"
$ws.Range("E41").Value = "  -1.02%  "

# Row 42
$ws.Range("E42").Value = "  -2.90%  "

# Row 43
$ws.Range("E43").Value = "  -1.74%  "

# Row 44
$ws.Range("D44").Value = "'0.0455"
$ws.Range("E44").Value = "  -0.44%  "

# Row 45
$ws.Range("E45").Value = "  +2.41%  "

# Row 46
$ws.Range("D46").Value = "'2.88"
$ws.Range("E46").Value = "  +0.97%  "

# Row 47
$ws.Range("E47").Value = "  +0.22%  "

# Row 48
$ws.Range("D48").Value = "'8.62"
$ws.Range("E48").Value = "  -0.80%  "

# Row 49
$ws.Range("E49").Value = "  -0.18%  "

# Row 50
$ws.Range("D50").Value = "'0.000251"
$ws.Range("E50").Value = "  -1.32%  "

# Row 51
$ws.Range("E51").Value = "  +0.18%  "
